$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order) that receive new values for rows 2..25
$cols = @("B","C","D","E","F","I","J","K","L","M")

# New values for each data row (row 2 through row 25), one array per row,
# values ordered to match $cols
$data = @(
    @(1.02,1.0620341478923,1.06547544941012,1.057539734874008,1.075063029343669,1.056633681462053,1.067006753425094,1.068188589065256,1.060274426405623,1.077750589779761),
    @(1.02,1.063255706570799,1.066448501941961,1.058595181580479,1.076178254674246,1.057076319578798,1.067881291922823,1.068976653676612,1.06114311578966,1.07868232810221),
    @(1.02,1.064045805260526,1.067077789068305,1.059277950703991,1.076899849272988,1.057361299230851,1.068446298854742,1.069485620964163,1.061704438951587,1.079284600288028),
    @(1.02,1.064377885438887,1.067342260232644,1.059564945989194,1.077203201433992,1.057480761396079,1.068683618790739,1.069699361220946,1.061940234318301,1.079537647025395),
    @(1.02,1.064433638687872,1.067386661374285,1.05961313139825,1.077254135266193,1.057500799509336,1.068723453632728,1.069735235716427,1.061979814588259,1.079580126010292),
    @(1.02,1.064050242833342,1.067081323262593,1.059281785708123,1.076903902702357,1.05736289683814,1.068449470754976,1.069488477872817,1.061707590385917,1.079287982093015),
    @(1.02,1.062447049226222,1.065804368649618,1.057896465057483,1.075439932032036,1.056783571004551,1.067302490317587,1.068455119234661,1.060568165595361,1.078065605137442),
    @(1.02,1.059619387104342,1.063551541878763,1.05545395713438,1.072859940843583,1.055751689355346,1.065274585850064,1.066626786943491,1.058554352769365,1.075906786186152),
    @(1.02,1.05773238043431,1.062047799567593,1.053824604729988,1.071139677541681,1.05505630348436,1.063918014069803,1.065402842981739,1.057207703999754,1.074464251830992),
    @(1.02,1.056914805722599,1.06139620709749,1.053118819169871,1.070394703023538,1.054753412690848,1.063329485991424,1.064871647354109,1.056623599350476,1.073838814068383),
    @(1.02,1.056611046047587,1.061154105837756,1.052856617534344,1.070117971298562,1.054640636544558,1.06311071000988,1.064674152980021,1.056406485574757,1.07360637528494),
    @(1.02,1.056676207011924,1.061206040605969,1.05291286257694,1.07017733189873,1.054664839600129,1.063157645940268,1.064716524572786,1.056453064088092,1.073656239799469),
    @(1.02,1.056889698419338,1.06137619638419,1.053097146355833,1.070371828608618,1.054744096072973,1.0633114053872,1.064855326182806,1.056605655753454,1.0738196031312),
    @(1.02,1.057021227368147,1.061481025551552,1.053210684133011,1.07049166241196,1.05479289292482,1.063406119021624,1.064940821923644,1.056699652419677,1.073920240280627),
    @(1.02,1.057786629659598,1.062091033760145,1.053871439668225,1.071189117029986,1.055076367676162,1.063957048927608,1.065438070866302,1.057246447968113,1.074505742861162),
    @(1.02,1.058266613672786,1.062473551099139,1.054285842018204,1.071626586869913,1.055253705549756,1.06430233073364,1.065749654421702,1.057589170401001,1.074872794991483),
    @(1.02,1.058546533286786,1.062696622497074,1.054527530418179,1.071881747135025,1.055356971655699,1.064503619476106,1.065931278239379,1.057788978521527,1.075086812156541),
    @(1.02,1.058641970759215,1.062772676573173,1.054609935643106,1.071968748799833,1.055392153571449,1.064572235416532,1.065993187362879,1.057857091676989,1.075159773272205),
    @(1.02,1.058215120806345,1.062432515221887,1.054241383216015,1.071579651387315,1.055234696680411,1.064265296477692,1.065716236629762,1.057552409477065,1.074833421865235),
    @(1.02,1.056826832640962,1.061326091714134,1.053042880533349,1.070314554602924,1.054720764451225,1.063266131780921,1.064814457677178,1.056560725500279,1.073771500119099),
    @(1.02,1.05595351959925,1.060630029005609,1.052289093692907,1.069519050340349,1.054396077914056,1.062636930294857,1.064246404337744,1.055936338276087,1.073103113341124),
    @(1.02,1.056416522071331,1.060999064372942,1.052688713706647,1.069940771058407,1.0545683482973,1.062970576013988,1.064547641970695,1.056267421320325,1.07345750604917),
    @(1.02,1.058238388361249,1.062451057706023,1.054261472313831,1.07160085953792,1.055243286501953,1.064282030998488,1.065731337054903,1.057569020453152,1.074851213127678),
    @(1.02,1.060350730625428,1.06413427442144,1.056085577139642,1.073526972291281,1.056019767902092,1.065799659064444,1.067100340975131,1.059075690213774,1.076465472874762)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowValues[$j]
    }
}
